# week1 DQ1 draft finished: update actual-time-to-complete for the
# "Discussion question 1" row and move the active selection to C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week1")

# C6 = "Actual time length to complete" for the "Discussion question 1" task.
# Previously 4.7222222222222221E-2 (1:08), now 8.0555555555555561E-2 (1:56).
$ws.Range("C6").Value = 0.08055555555555556

# C20 is SUM(C2:C19) and recalculates automatically from the C6 change.

# Move the active cell/selection on the sheet from C13 to C7.
$ws.Range("C7").Select()
